$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New table contents (header + 17 player rows) replacing the previous
# 18-row roster. Two players (Clint Capela, Bogdan Bogdanovic) were
# dropped and a new player (Alexandre Sarr) was added, and the whole
# roster was re-ordered.
$data = @(
    @("Oyuncu Adı", "Pozisyon", "Takım"),
    @("Dejounte Murray", "PG,SG", "New Orleans Pelicans"),
    @("Jalen Suggs", "PG,SG", "Orlando Magic"),
    @("Chris Paul", "PG", "San Antonio Spurs"),
    @("Jaylen Brown", "SG,SF", "Boston Celtics"),
    @("Deni Avdija", "SF,PF", "Portland Trail Blazers"),
    @("Pascal Siakam", "SF,PF", "Indiana Pacers"),
    @("Naz Reid", "PF,C", "Minnesota Timberwolves"),
    @("Jerami Grant", "SF,PF", "Portland Trail Blazers"),
    @("Rudy Gobert", "C", "Minnesota Timberwolves"),
    @("Alexandre Sarr", "PF,C", "Washington Wizards"),
    @("Jakob Poeltl", "C", "Toronto Raptors"),
    @("Nikola Jokic", "C", "Denver Nuggets"),
    @("Russell Westbrook", "PG", "Denver Nuggets"),
    @("Ayo Dosunmu", "SG,SF", "Chicago Bulls"),
    @("Jalen Green", "PG,SG", "Houston Rockets"),
    @("Paolo Banchero", "SF,PF", "Orlando Magic"),
    @("Chet Holmgren", "PF,C", "Oklahoma City Thunder")
)

$rowCount = $data.Length

# Remove the old row 19 (one fewer data row than before: 18 -> 17).
$ws.Rows.Item(19).Delete() | Out-Null

for ($i = 0; $i -lt $rowCount; $i++) {
    $r = $i + 1
    $ws.Cells.Item($r, 1).Value = $data[$i][0]
    $ws.Cells.Item($r, 2).Value = $data[$i][1]
    $ws.Cells.Item($r, 3).Value = $data[$i][2]
}
